$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("M2").Value = 0.5273236666666667
$ws.Range("N2").Value = 1.581971
$ws.Range("O2").Value = 0.1756607078989806
$ws.Range("P2").Value = 0.1756607078989806
$ws.Range("Q2").Value = 3.272578760962889
$ws.Range("R2").Value = 29.45320884866601
$ws.Range("S2").Value = 0.002637939082102384
$ws.Range("T2").Value = 0.002637939082102384
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("O3").Value = 0.4904710091781626
$ws.Range("P3").Value = 0.4904710091781627
$ws.Range("Q3").Value = 9.137530109622222
$ws.Range("R3").Value = 82.2377709866
$ws.Range("S3").Value = 0.007365521061735289
$ws.Range("T3").Value = 0.007365521061735289
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 1.002254
$ws.Range("N4").Value = 3.006762
$ws.Range("O4").Value = 0.3338682829228568
$ws.Range("P4").Value = 0.3338682829228568
$ws.Range("Q4").Value = 6.220003691894666
$ws.Range("R4").Value = 55.980033227052
$ws.Range("S4").Value = 0.005013780271813027
$ws.Range("T4").Value = 0.005013780271813027
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("M5").Value = 0.5273236666666667
$ws.Range("N5").Value = 1.581971
$ws.Range("O5").Value = 0.1756607078989806
$ws.Range("P5").Value = 0.1756607078989806
$ws.Range("Q5").Value = 203.0472805152313
$ws.Range("R5").Value = 1827.425524637081
$ws.Range("S5").Value = 0.1636710361794737
$ws.Range("T5").Value = 0.1636710361794737
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("O6").Value = 0.4904710091781626
$ws.Range("P6").Value = 0.4904710091781627
$ws.Range("S6").Value = 0.4569940497697829
$ws.Range("T6").Value = 0.4569940497697829
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 1.002254
$ws.Range("N7").Value = 3.006762
$ws.Range("O7").Value = 0.3338682829228568
$ws.Range("P7").Value = 0.3338682829228568
$ws.Range("Q7").Value = 385.9203786014647
$ws.Range("R7").Value = 3473.283407413182
$ws.Range("S7").Value = 0.3110801981105005
$ws.Range("T7").Value = 0.3110801981105006
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("M8").Value = 0.5273236666666667
$ws.Range("N8").Value = 1.581971
$ws.Range("O8").Value = 0.1756607078989806
$ws.Range("P8").Value = 0.1756607078989806
$ws.Range("Q8").Value = 11.60158769966078
$ws.Range("R8").Value = 104.414289296947
$ws.Range("S8").Value = 0.009351732637404509
$ws.Range("T8").Value = 0.009351732637404509
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("O9").Value = 0.4904710091781626
$ws.Range("P9").Value = 0.4904710091781627
$ws.Range("Q9").Value = 32.39337069274444
$ws.Range("R9").Value = 291.5403362347
$ws.Range("S9").Value = 0.02611143834664444
$ws.Range("T9").Value = 0.02611143834664444
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 1.002254
$ws.Range("N10").Value = 3.006762
$ws.Range("O10").Value = 0.3338682829228568
$ws.Range("P10").Value = 0.3338682829228568
$ws.Range("Q10").Value = 22.05047566295933
$ws.Range("R10").Value = 198.454280966634
$ws.Range("S10").Value = 0.01777430454054319
$ws.Range("T10").Value = 0.01777430454054319
